$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0864294738795979, 0.1717255535034242, 0.05132172044861762, 0.2265429770454551, 0.2173128427139068)
    3  = @(0.1015794586745889, 0.1777261608473147, 0.06528206036901001, 0.255503542771935, 0.244016352129155)
    4  = @(0.09564297872977319, 0.182775340538003, 0.07239544141620145, 0.2690640098864979, 0.2626740711899914)
    5  = @(0.0515887520981297, 0.1373511933772956, 0.03807726663693461, 0.1951339709966838, 0.1973764272229122)
    6  = @(0.03412360743610758, 0.1315810182071815, 0.02988307038476774, 0.1728672044801088, 0.1786326142932414)
    7  = @(0.03856004825347423, 0.1426784820072856, 0.04033731687300466, 0.2008415217852241, 0.2090615806303401)
    8  = @(0.06589614311488577, 0.08267801794971415, 0.01154809474451304, 0.107462061884709, 0.09298898687755791)
    9  = @(-0.03841838264255042, 0.03841838264255042, 0.001596828758587614, 0.03996033982072242, 0.01346420998712108)
    10 = @(-0.02233679963354765, 0.02233679963354765, 0.0004989326178692543, 0.02233679963354765)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    if ($values.Count -ge 5) {
        $ws.Cells.Item($row, 6).Value = $values[4]
    }
}
